$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44784
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("Q2").Value = "`$/bandeja 18 kilos"
$ws.Range("S2").Value = 1083
$ws.Range("T2").Value = 18
$ws.Range("D3").Value = 45043
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21500
$ws.Range("S3").Value = 1194
$ws.Range("D4").Value = 44616
$ws.Range("L4").Value = "Segunda"
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 16500
$ws.Range("Q4").Value = "`$/caja 18 kilos granel"
$ws.Range("S4").Value = 917
$ws.Range("T4").Value = 18
$ws.Range("D5").Value = 45002
$ws.Range("L5").Value = "Segunda"
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 24500
$ws.Range("S5").Value = 1361
$ws.Range("D6").Value = 45086
$ws.Range("L6").Value = "Especial"
$ws.Range("P6").Value = 25500
$ws.Range("S6").Value = 1417
$ws.Range("D7").Value = 45086
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 21000
$ws.Range("P7").Value = 20500
$ws.Range("S7").Value = 1139
$ws.Range("D8").Value = 45034
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 25000
$ws.Range("O8").Value = 26000
$ws.Range("P8").Value = 25600
$ws.Range("S8").Value = 1422
$ws.Range("D9").Value = 44491
$ws.Range("L9").Value = "Primera"
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("Q9").Value = "`$/bandeja 10 kilos"
$ws.Range("S9").Value = 1450
$ws.Range("T9").Value = 10
$ws.Range("D10").Value = 44819
$ws.Range("L10").Value = "Primera"
$ws.Range("Q10").Value = "`$/bandeja 10 kilos"
$ws.Range("S10").Value = 1750
$ws.Range("T10").Value = 10
$ws.Range("D11").Value = 44656
$ws.Range("M11").Value = 270
$ws.Range("N11").Value = 19000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 19500
$ws.Range("Q11").Value = "`$/bandeja 18 kilos"
$ws.Range("S11").Value = 1083
$ws.Range("T11").Value = 18
$ws.Range("D12").Value = 44614
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 21000
$ws.Range("P12").Value = 20500
$ws.Range("Q12").Value = "`$/bandeja 18 kilos"
$ws.Range("S12").Value = 1139
$ws.Range("D13").Value = 44489
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 26000
$ws.Range("O13").Value = 27000
$ws.Range("P13").Value = 26500
$ws.Range("S13").Value = 1472
$ws.Range("D14").Value = 44263
$ws.Range("N14").Value = 21000
$ws.Range("O14").Value = 22000
$ws.Range("P14").Value = 21500
$ws.Range("Q14").Value = "`$/caja 18 kilos"
$ws.Range("S14").Value = 1194
$ws.Range("D15").Value = 44487
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 14500
$ws.Range("Q15").Value = "`$/bandeja 10 kilos"
$ws.Range("S15").Value = 1450
$ws.Range("T15").Value = 10
$ws.Range("D16").Value = 44418
$ws.Range("M16").Value = 240
$ws.Range("N16").Value = 10000
$ws.Range("O16").Value = 11000
$ws.Range("P16").Value = 10500
$ws.Range("S16").Value = 1050
$ws.Range("D17").Value = 44789
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 250
$ws.Range("D18").Value = 44673
$ws.Range("L18").Value = "Especial"
$ws.Range("M18").Value = 400
$ws.Range("N18").Value = 14000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 14500
$ws.Range("Q18").Value = "`$/bandeja 10 kilos"
$ws.Range("S18").Value = 1450
$ws.Range("T18").Value = 10
$ws.Range("D19").Value = 44706
$ws.Range("M19").Value = 400
$ws.Range("N19").Value = 9000
$ws.Range("O19").Value = 10000
$ws.Range("P19").Value = 9500
$ws.Range("S19").Value = 950
$ws.Range("D20").Value = 44629
$ws.Range("N20").Value = 17000
$ws.Range("O20").Value = 18000
$ws.Range("P20").Value = 17500
$ws.Range("S20").Value = 972
$ws.Range("D21").Value = 45069
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("M21").Value = 370
$ws.Range("P21").Value = 19486
$ws.Range("R21").Value = "Región Metropolitana"
$ws.Range("D22").Value = 44602
$ws.Range("M22").Value = 270
$ws.Range("N22").Value = 20000
$ws.Range("O22").Value = 21000
$ws.Range("P22").Value = 20500
$ws.Range("Q22").Value = "`$/bandeja 18 kilos"
$ws.Range("S22").Value = 1139
$ws.Range("T22").Value = 18
$ws.Range("D23").Value = 44323
$ws.Range("M23").Value = 270
$ws.Range("N23").Value = 21000
$ws.Range("O23").Value = 22000
$ws.Range("P23").Value = 21500
$ws.Range("S23").Value = 1194
$ws.Range("D24").Value = 44307
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 250
$ws.Range("N24").Value = 19000
$ws.Range("O24").Value = 20000
$ws.Range("P24").Value = 19500
$ws.Range("Q24").Value = "`$/bandeja 18 kilos"
$ws.Range("S24").Value = 1083
$ws.Range("D25").Value = 44291
$ws.Range("K25").Value = "Hayward"
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 17000
$ws.Range("O25").Value = 18000
$ws.Range("P25").Value = 17500
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 972
$ws.Range("D26").Value = 44991
$ws.Range("N26").Value = 24000
$ws.Range("O26").Value = 25000
$ws.Range("P26").Value = 24500
$ws.Range("S26").Value = 1361
